# Re-apply the refreshed cryptos snapshot (prices / 1h volume deltas) from the
# GitHub Actions data pull onto Sheet1. Column A (rank) / B (coin) / C (link) are
# mostly unchanged except for the Aave <-> TrustWalletToken row swap below.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.619.19'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '1.595.98'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').Value = "'" + '210.99'
$ws.Range('E5').Value = '  -0.49%  '
$ws.Range('D6').Value = "'" + '0.514'
$ws.Range('E6').Value = '  -0.13%  '
$ws.Range('E7').Value = '  +0.15%  '
$ws.Range('E8').Value = '  +0.10%  '
$ws.Range('E9').Value = '  -0.30%  '
$ws.Range('D10').Value = "'" + '19.48'
$ws.Range('E10').Value = '  -1.30%  '
$ws.Range('E11').Value = '  +0.00%  '
$ws.Range('D12').Value = '1.820.40'
$ws.Range('D13').Value = '1.576.02'
$ws.Range('E13').Value = '  -1.25%  '
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('D16').Value = "'" + '65.01'
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('D17').Value = '26.603.56'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').Value = '0.0₃0739'
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').Value = "'" + '208.44'
$ws.Range('E20').Value = '  -0.89%  '
$ws.Range('D21').Value = "'" + '7.02'
$ws.Range('E21').Value = '  +4.79%  '
$ws.Range('D23').Value = "'" + '2.29'
$ws.Range('E23').Value = '  -1.08%  '
$ws.Range('D24').Value = "'" + '8.90'
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('D25').Value = "'" + '145.13'
$ws.Range('E25').Value = '  -1.01%  '
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('D27').Value = "'" + '7.10'
$ws.Range('E27').Value = '  -1.28%  '
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('D29').Value = "'" + '15.25'
$ws.Range('E29').Value = '  -0.57%  '
$ws.Range('E30').Value = '  +0.82%  '
$ws.Range('D31').Value = "'" + '1.15'
$ws.Range('E31').Value = '  -0.21%  '
$ws.Range('E32').Value = '  -0.48%  '
$ws.Range('E33').Value = '  +0.35%  '
$ws.Range('D34').Value = '1.276.35'
$ws.Range('E34').Value = '  -1.17%  '
$ws.Range('E35').Value = '  -8.19%  '
$ws.Range('E36').Value = '  +0.89%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('E38').Value = '  -0.64%  '
$ws.Range('D39').Value = "'" + '0.838'
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('D40').Value = "'" + '1.04'
$ws.Range('E40').Value = '  +18.84%  '
$ws.Range('D41').Value = "'" + '5.48'
$ws.Range('E41').Value = '  +1.67%  '
$ws.Range('E42').Value = '  +0.27%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = "'" + '0.784'
$ws.Range('E43').Value = '  -1.00%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = "'" + '64.16'
$ws.Range('E44').Value = '  +0.51%  '
$ws.Range('D45').Value = '1.732.61'
$ws.Range('D46').Value = "'" + '90.12'
$ws.Range('E46').Value = '  +0.31%  '
$ws.Range('E48').Value = '  +3.66%  '
$ws.Range('E49').Value = '  +0.89%  '
$ws.Range('E50').Value = '  +0.17%  '
$ws.Range('D51').Value = "'" + '7.43'
$ws.Range('E51').Value = '  -1.48%  '
